$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 7 (USB Micro-B Power Breakout): source changes from "ebay" to "Adafruit",
# and the cost figures are updated.
# ---------------------------------------------------------------------------
$ws.Range("C7").Value = "Adafruit"
$ws.Range("D7").Value = 33.75
$ws.Range("E7").Value = 4.8
$ws.Range("F7").Value = 25
$ws.Range("G7").Formula = "=(D7+E7)/F7"

# ---------------------------------------------------------------------------
# Row 6 (Arduino Pro Mini) gets a highlight fill across the whole row.
# ---------------------------------------------------------------------------
$ws.Range("A6:G6").Interior.ThemeColor = 6

# ---------------------------------------------------------------------------
# Row 8: new line item "SD card breakout" (was a lone label) with full detail,
# highlighted with the alternate (white) fill.
# ---------------------------------------------------------------------------
$ws.Range("B8").Value = "DEV-13743"
$ws.Range("C8").Value = "SparkFun"
$ws.Range("D8").Value = 111.5
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 25
$ws.Range("G8").Formula = "=(D8+E8)/F8"
$ws.Range("H8").Value = "*"
$ws.Range("A8:G8").Interior.ThemeColor = 2

# ---------------------------------------------------------------------------
# Row 9: new line item "Micro SD card (2GB)" (was a lone label) with detail.
# ---------------------------------------------------------------------------
$ws.Range("B9").Value = "GND2353"
$ws.Range("C9").Value = "Unique Photo"
$ws.Range("D9").Value = 174.75
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 25
$ws.Range("G9").Formula = "=(D9+E9)/F9"
$ws.Range("H9").Value = "*"

# ---------------------------------------------------------------------------
# Row 10: "uBlox GPS" (was a lone label) with detail + highlight fill.
# ---------------------------------------------------------------------------
$ws.Range("D10").Value = 312.5
$ws.Range("E10").Value = 13.57
$ws.Range("F10").Value = 25
$ws.Range("G10").Formula = "=(D10+E10)/F10"
$ws.Range("A10:G10").Interior.ThemeColor = 6

# ---------------------------------------------------------------------------
# Row 11: new line item "Nylon tactical belt".
# ---------------------------------------------------------------------------
$ws.Range("A11").Value = "Nylon tactical belt"
$ws.Range("C11").Value = "Amazon"
$ws.Range("D11").Value = 7.99
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 1
$ws.Range("G11").Formula = "=(D11+E11)/F11"
$ws.Range("H11").Value = "*"

# ---------------------------------------------------------------------------
# Row 12: new line item "Waterproof airtight survival case".
# ---------------------------------------------------------------------------
$ws.Range("A12").Value = "Waterproof airtight survival case"
$ws.Range("C12").Value = "Amazon"
$ws.Range("D12").Value = 3.28
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 1
$ws.Range("G12").Formula = "=(D12+E12)/F12"

# ---------------------------------------------------------------------------
# Row 13: new line item "5v 5600 mAh USB battery".
# ---------------------------------------------------------------------------
$ws.Range("A13").Value = "5v 5600 mAh USB battery"
$ws.Range("C13").Value = "Amazon"
$ws.Range("D13").Value = 10
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 1
$ws.Range("G13").Formula = "=(D13+E13)/F13"

# ---------------------------------------------------------------------------
# Rows 14-16: blank placeholder rows that still compute into the total.
# ---------------------------------------------------------------------------
$ws.Range("F14").Value = 1
$ws.Range("G14").Formula = "=(D14+E14)/F14"
$ws.Range("F15").Value = 1
$ws.Range("G15").Formula = "=(D15+E15)/F15"
$ws.Range("F16").Value = 1
$ws.Range("G16").Formula = "=(D16+E16)/F16"

# ---------------------------------------------------------------------------
# Row 4 (TP5110 Low Power Timer Breakout) picks up the "*" marker too.
# ---------------------------------------------------------------------------
$ws.Range("H4").Value = "*"

# ---------------------------------------------------------------------------
# Row 17: TOTAL row (sum picks up the new rows automatically via SUM(G2:G16)).
# ---------------------------------------------------------------------------
$ws.Range("G17").Formula = "=SUM(G2:G16)"

$wb.Application.Calculate()

# ---------------------------------------------------------------------------
# Selection / active cell, matching the saved view state.
# ---------------------------------------------------------------------------
$ws.Range("E18").Select()
